# Applies the documented diff:
#  1. Strip the bare "_GoBack" bookmark paragraph (after "Modules: pyqtdarktheme")
#     down to a totally empty paragraph.
#  2. Rework the "Morgen:" .. "Taktlogik (...)" block at the end of the document
#     into the new "Speichern & Öffnen" / "Konzept:" / two narrative paragraphs,
#     with the "_GoBack" bookmark now trailing the very last run.

$d = $word.ActiveDocument

# --- Change 1 : drop the orphan _GoBack bookmark paragraph right after
#     "Modules: pyqtdarktheme" -- becomes a plain empty <w:p/>.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- Change 2 : rebuild the tail of the document -------------------------

# Locate the three paragraphs we are rewriting: "Morgen:", "Seitenzahl",
# "Taktlogik (...)" by searching from the back of the document (robust to
# any paragraph-index drift caused by Change 1).
$total = $d.Paragraphs.Count
$pTaktlogik = $d.Paragraphs.Item($total - 1)
$pSeitenzahl = $d.Paragraphs.Item($total - 2)
$pMorgen = $d.Paragraphs.Item($total - 3)

# Drop "Seitenzahl" and "Taktlogik (...)" completely -- their whole ranges
# (including paragraph marks) disappear, leaving "Morgen:" followed directly
# by the trailing blank paragraph.
$killRange = $d.Range($pSeitenzahl.Range.Start, $pTaktlogik.Range.End)
$killRange.Delete()

# "Morgen:" -> "Speichern & Öffnen" (keeps the paragraph's en-GB run formatting).
$d.Content.Find.Execute("Morgen:", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Speichern & Öffnen", 2)

# New "Konzept:" paragraph, inserted right after "Speichern & Öffnen" --
# inherits that paragraph's (en-GB) formatting, which is what the target needs.
$pSpeichern = $d.Paragraphs.Item($total - 3)
$pSpeichern.Range.InsertParagraphAfter()
$pKonzept = $d.Paragraphs.Item($total - 2)
$pKonzept.Range.Text = "Konzept:"

# Two new plain-formatted narrative paragraphs, inserted right before the
# document's trailing blank paragraph so they pick up no special formatting.
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$pLast.Range.InsertParagraphBefore()
$pLast.Range.InsertParagraphBefore()

$pPara3 = $d.Paragraphs.Item($d.Paragraphs.Count - 2)
$pPara3.Range.Text = "Da das Programm sehr komplexe und vernetzte Klassen beinhaltet, ist die Entscheidung getroffen worden, nicht, wie anfangs beabsichtigt, alle member einer klasse in ein dict zu schreiben, welches in eine json-Datei geschrieben werden hätte können, sonder:"

$pPara4 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$pPara4.Range.Text = "Die Bibliothek pickle zu verwenden, mit der es möglich ist, Objekte von Klassen zu speichern und zu laden, wärend alle member erhalten bleiben (sowohl deren Bezeichnung als auch der Inhalt). Einige Klasse erben aber von QGraphics__ - Klassen, deren Objekte nicht von Pickle unterstützt werden. Daher werden diese Vererbungen ausgelagert, indem alle zusätzlichen objekte einer klasse in ein Index-based dict verlagert werden und ein index zurückgegeben wird, welcher später dazu verwendet werden kann, das eigentliche graphics item wieder zu bekommen."

# The _GoBack bookmark now sits at the very end of that last paragraph's text.
$d.Bookmarks.Add("_GoBack", $pPara4.Range)
